# Update the "Förändrad" date column (C) for rows 2-39 from 45188 to 45189
# (i.e. bump the serial date by one day, 2023-09-19 -> 2023-09-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 39; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45188) {
        $cell.Value = 45189
    }
}
